$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Producto"
$ws.Range("B1").Value = "Precio"
$ws.Range("C1").Value = "Cantidad"

# Data row -- Precio/Cantidad are numeric-looking but stored as text in the
# source workbook, so prefix with an apostrophe to force text storage
# instead of Excel auto-converting them to numbers.
$ws.Range("A2").Value = "Camisa"
$ws.Range("B2").Value = "'5000"
$ws.Range("C2").Value = "'10"
